# Weekly CompStat (114th Precinct) data refresh: new crime data collected.
# Updates the report header (volume number + week-covering dates) and the
# full crime-complaints table (rows 15-28) with the new week's figures,
# including a handful of cells that flip between the literal placeholder
# text ("0" / "***.*") and real numbers as data becomes available/unavailable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 32   Number  9" -> "...Number  10"
#         "Report Covering the Week  2/24/2025  Through  3/2/2025"
#         -> "...3/3/2025  Through  3/9/2025"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/3/2025  Through  3/9/2025"

# ---------------------------------------------------------------------
# Helper cell used to author literal TEXT values ("0" / "***.*") into
# cells without Excel's autodetection turning numeric-looking text back
# into a number. We build the text with a formula (forces text result),
# paste-special VALUES into the target (keeps it text), then paste-special
# FORMATS from a cell that already carries the desired display style.
# ---------------------------------------------------------------------
function Set-TextValue($cellRef, $text, $formatSourceRef) {
    $helper = $ws.Range("ZZ1")
    $helper.Formula = '=T("' + $text + '")'
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
    $ws.Range($formatSourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
    $helper.Clear()
}

function Set-NumberValue($cellRef, $value, $formatSourceRef) {
    $ws.Range($formatSourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($cellRef).Value = $value
}

# ---------------------------------------------------------------------
# Row 15 (Rape): C/D go from "0" text to real numbers; E goes from
# "***.*" text to a real percentage (0).
# ---------------------------------------------------------------------
Set-NumberValue "C15" 1 "F15"
Set-NumberValue "D15" 1 "F15"
Set-NumberValue "E15" 0 "H15"

# Row 16 (Robbery): C goes from "0" text to a real number.
Set-NumberValue "C16" 2 "D16"

# Row 22 (Transit): C goes from "0" text to a real number; D/E flip the
# other way, from real numbers to the "0" / "***.*" placeholder text.
Set-NumberValue "C22" 1 "F22"
Set-TextValue "D22" "0" "C14"
Set-TextValue "E22" "***.*" "E14"

# Row 27 (UCR Rape*): C/D go from "0" text to real numbers; E goes from
# "***.*" text to a real percentage (0).
Set-NumberValue "C27" 1 "F27"
Set-NumberValue "D27" 1 "F27"
Set-NumberValue "E27" 0 "H27"

# Rows 29-31 (Shooting Vic. / Shooting Inc. / Hate Crimes): D/E flip from
# real numbers to the "0" / "***.*" placeholder text.
Set-TextValue "D29" "0" "C14"
Set-TextValue "E29" "***.*" "E14"
Set-TextValue "D30" "0" "C14"
Set-TextValue "E30" "***.*" "E14"
Set-TextValue "D31" "0" "C14"
Set-TextValue "E31" "***.*" "E14"

# ---------------------------------------------------------------------
# Remaining numeric-only updates (style unchanged) across rows 15-28.
# ---------------------------------------------------------------------
# Row 15
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = -14.285714285714
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -25
# Row 16
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = -80.555555555555
$ws.Range("I16").Value = 24
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = -60
$ws.Range("L16").Value = -45.454545454545
$ws.Range("M16").Value = -59.322033898305
$ws.Range("N16").Value = -93.684210526315
# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 51
$ws.Range("G17").Value = 66
$ws.Range("H17").Value = -22.727272727272
$ws.Range("I17").Value = 115
$ws.Range("J17").Value = 112
$ws.Range("K17").Value = 2.678571428571
$ws.Range("L17").Value = 59.722222222222
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = -22.297297297297
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -22.727272727272
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 48
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = -30.434782608695
$ws.Range("M18").Value = -72.413793103448
$ws.Range("N18").Value = -92.270531400966
# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -31.578947368421
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 89
$ws.Range("H19").Value = -37.078651685393
$ws.Range("I19").Value = 154
$ws.Range("J19").Value = 176
$ws.Range("K19").Value = -12.5
$ws.Range("L19").Value = -8.333333333333
$ws.Range("M19").Value = 49.514563106796
$ws.Range("N19").Value = -27.699530516431
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -41.666666666666
$ws.Range("I20").Value = 41
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = -26.785714285714
$ws.Range("L20").Value = -28.070175438596
$ws.Range("M20").Value = 10.810810810810
$ws.Range("N20").Value = -93.086003372681
# Row 21 (TOTAL, bold)
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -31.707317073170
$ws.Range("F21").Value = 146
$ws.Range("G21").Value = 241
$ws.Range("H21").Value = -39.419087136929
$ws.Range("I21").Value = 372
$ws.Range("J21").Value = 459
$ws.Range("K21").Value = -18.954248366013
$ws.Range("L21").Value = -5.583756345177
$ws.Range("M21").Value = 1.086956521739
$ws.Range("N21").Value = -78.875638841567
# Row 22
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 75
$ws.Range("L22").Value = 133.333333333333
$ws.Range("M22").Value = 16.666666666666
# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = -35.714285714285
$ws.Range("I23").Value = 42
$ws.Range("J23").Value = 41
$ws.Range("K23").Value = 2.439024390243
$ws.Range("L23").Value = 7.692307692307
$ws.Range("M23").Value = 55.555555555555
# Row 24
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = 14.583333333333
$ws.Range("F24").Value = 192
$ws.Range("G24").Value = 177
$ws.Range("H24").Value = 8.474576271186
$ws.Range("I24").Value = 435
$ws.Range("J24").Value = 428
$ws.Range("K24").Value = 1.635514018691
$ws.Range("L24").Value = 1.873536299765
$ws.Range("M24").Value = 97.727272727272
# Row 25
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 29
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 112
$ws.Range("G25").Value = 92
$ws.Range("H25").Value = 21.739130434782
$ws.Range("I25").Value = 245
$ws.Range("J25").Value = 211
$ws.Range("K25").Value = 16.113744075829
$ws.Range("L25").Value = 16.113744075829
# Row 26
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 157.142857142857
$ws.Range("F26").Value = 68
$ws.Range("G26").Value = 66
$ws.Range("H26").Value = 3.030303030303
$ws.Range("I26").Value = 161
$ws.Range("J26").Value = 156
$ws.Range("K26").Value = 3.205128205128
$ws.Range("L26").Value = 13.380281690140
$ws.Range("M26").Value = -0.617283950617
# Row 27
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -22.222222222222
$ws.Range("L27").Value = -22.222222222222
# Row 28
$ws.Range("C28").Value = 3
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 60
$ws.Range("L28").Value = 23.076923076923
# Row 29
$ws.Range("L29").Value = 0
# Row 30
$ws.Range("L30").Value = 0
